# Applies the two changes captured by the commit's XML diff:
#
#   1. The table on slide 16 switches from table style
#      {FA496728-B941-4398-B0E9-E2A6313EB5A9} ("Table_0" built-in style)
#      to {34945F98-ED17-4964-93EA-C8C16287FD5F}.
#
#   2. The deck's "Integral" and "Office Theme" color schemes are swapped
#      between the two theme parts (theme1.xml <-> theme2.xml). The part
#      that actually drives the rendered deck (the one behind the slide
#      master / presentation relationship) picks up the "Office Theme"
#      palette; here we reproduce that by rewriting its 12 theme colors.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style swap on slide 16
# ---------------------------------------------------------------------
$targetStyleId = "{34945F98-ED17-4964-93EA-C8C16287FD5F}"
$slide = $p.Slides.Item(16)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($targetStyleId)
    }
}

# ---------------------------------------------------------------------
# 2) Theme color scheme swap (Integral -> Office Theme)
# ---------------------------------------------------------------------
function HexToRgbValue($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# clrScheme slot order exposed by ThemeColorScheme.Item(1..12):
#   dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$scheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $scheme.Count; $i++) {
    $scheme.Item($i).RGB = HexToRgbValue $officeThemeColors[$i - 1]
}
